$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E (Comments) ---
$ws.Columns.Item(5).ColumnWidth = 54

# Header E1 "Коментарі" - same look as the other headers (22pt font, border, centered-ish) but a new accent colour
$ws.Range("E1").Value = "Коментарі"

# Data E2 "Трабл з підключенням файлів ресурсів"
$ws.Range("E2").Value = "Трабл з підключенням файлів ресурсів"

# --- Row 2 gets an owner + status for the first task ---
$ws.Range("B2").Value = "Оля"
$ws.Range("C2").Value = "process"

# --- Formatting ---
# Reuse column B / column C existing per-column body styles for B2 / C2
$ws.Range("B1").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$ws.Range("B2").Value = "Оля"

$ws.Range("C1").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$ws.Range("C2").Value = "process"

# Give the new Comments column its own accent colour (like the other columns)
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$ws.Range("E1").Value = "Коментарі"
$ws.Range("E1").Interior.ThemeColor = 8
$ws.Range("E1").Interior.TintAndShade = 0.79998168889431442

$ws.Range("D2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null
$ws.Range("E2").Value = "Трабл з підключенням файлів ресурсів"
$ws.Range("E2").Interior.ThemeColor = 8
$ws.Range("E2").Interior.TintAndShade = 0.79998168889431442

$ws.Range("E1").Select()
$excel.CutCopyMode = $false
